$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '43.678.49'
$ws.Cells.Item(2, 5).Value = '  +2.84%  '

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '2.199.10'

$ws.Cells.Item(4, 5).Value = '  +0.15%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '258.66'

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '83.19'
$ws.Cells.Item(6, 5).Value = '  +10.72%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.620'
$ws.Cells.Item(7, 5).Value = '  +1.12%  '

$ws.Cells.Item(8, 5).Value = '  +0.04%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.598'
$ws.Cells.Item(9, 5).Value = '  +2.77%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '44.41'
$ws.Cells.Item(10, 5).Value = '  +10.36%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0920'
$ws.Cells.Item(11, 5).Value = '  +1.25%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '7.17'
$ws.Cells.Item(12, 5).Value = '  +5.99%  '

$ws.Cells.Item(13, 5).Value = '  +2.29%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '2.529.17'
$ws.Cells.Item(14, 5).Value = '  +0.62%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '14.38'
$ws.Cells.Item(15, 5).Value = '  +1.69%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '2.183.41'
$ws.Cells.Item(16, 5).Value = '  -0.15%  '

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.784'
$ws.Cells.Item(17, 5).Value = '  +2.47%  '

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '43.611.77'
$ws.Cells.Item(18, 5).Value = '  +2.89%  '

$ws.Cells.Item(19, 5).Value = '  +1.72%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '69.70'
$ws.Cells.Item(20, 5).Value = '  -1.39%  '

$ws.Cells.Item(21, 5).Value = '  +0.91%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '2.35'
$ws.Cells.Item(22, 5).Value = '  +11.17%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '231.54'
$ws.Cells.Item(23, 5).Value = '  +2.12%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '8.98'
$ws.Cells.Item(24, 5).Value = '  -4.64%  '

$ws.Cells.Item(25, 5).Value = '  -0.06%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '10.65'
$ws.Cells.Item(26, 5).Value = '  +2.15%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '3.45'
$ws.Cells.Item(27, 5).Value = '  +1.85%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '39.19'
$ws.Cells.Item(28, 5).Value = '  +3.61%  '

$ws.Cells.Item(29, 5).Value = '  -0.73%  '

$ws.Cells.Item(30, 5).Value = '  +3.39%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '174.42'
$ws.Cells.Item(31, 5).Value = '  +1.18%  '

$ws.Cells.Item(32, 5).Value = '  +1.98%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.0858'
$ws.Cells.Item(33, 5).Value = '  +4.61%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '5.32'
$ws.Cells.Item(34, 5).Value = '  +3.81%  '

$ws.Cells.Item(35, 5).Value = '  +2.16%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.111'
$ws.Cells.Item(36, 5).Value = '  +3.74%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '4.51'
$ws.Cells.Item(37, 5).Value = '  +7.09%  '

$ws.Cells.Item(38, 5).Value = '  +6.82%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '12.43'
$ws.Cells.Item(39, 5).Value = '  +4.54%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '2.82'
$ws.Cells.Item(40, 5).Value = '  +8.92%  '

$ws.Cells.Item(41, 5).Value = '  +1.69%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '62.95'
$ws.Cells.Item(42, 5).Value = '  +7.14%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '5.48'
$ws.Cells.Item(43, 5).Value = '  +6.27%  '

$ws.Cells.Item(44, 5).Value = '  +3.28%  '

$ws.Cells.Item(45, 2).Value = 'Cronos'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.0979'
$ws.Cells.Item(45, 5).Value = '  +0.87%  '

$ws.Cells.Item(46, 2).Value = 'FraxShare'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '8.30'
$ws.Cells.Item(46, 5).Value = '  +1.69%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '99.67'
$ws.Cells.Item(47, 5).Value = '  -1.54%  '

$ws.Cells.Item(48, 5).Value = '  +5.69%  '

$ws.Cells.Item(49, 5).Value = '  +1.89%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.441'
$ws.Cells.Item(50, 5).Value = '  -3.82%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '1.47'
$ws.Cells.Item(51, 5).Value = '  +7.74%  '
